# Update Name of Algo
# Applies updated imputed values produced by the RandomForest algorithm run
# to the terrestrial_mammals / combination_2_ABCDE / AD / 20 / seed3 result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = -7.1267
$ws.Range("A8").Value = -22.48010000000001
$ws.Range("A10").Value = -21.80439999999999
$ws.Range("A12").Value = -21.53779999999999
$ws.Range("D14").Value = -7.911400000000002
$ws.Range("D15").Value = -7.885999999999997
$ws.Range("A18").Value = -22.41590000000001
$ws.Range("D18").Value = -7.957499999999995
$ws.Range("D20").Value = -7.5337
$ws.Range("A25").Value = -21.55989999999998
$ws.Range("D29").Value = -7.080900000000002
$ws.Range("D30").Value = -7.124500000000005
$ws.Range("D31").Value = -8.452100000000002
$ws.Range("D35").Value = -7.913499999999996
$ws.Range("A37").Value = -19.1834
$ws.Range("D40").Value = -7.477899999999997
$ws.Range("D44").Value = -7.372899999999999
$ws.Range("D50").Value = -8.156399999999996
$ws.Range("D54").Value = -8.0274
$ws.Range("A55").Value = -22.1749
$ws.Range("A68").Value = -21.56509999999998
$ws.Range("D68").Value = -6.937799999999997
$ws.Range("D76").Value = -7.334200000000001
$ws.Range("A77").Value = -21.46399999999999
$ws.Range("A78").Value = -20.98079999999998
$ws.Range("A79").Value = -21.47829999999999
$ws.Range("A80").Value = -20.2519
$ws.Range("A81").Value = -21.7771
$ws.Range("A82").Value = -22.20830000000001
$ws.Range("A84").Value = -22.17340000000001
$ws.Range("D87").Value = -7.870599999999999
$ws.Range("D88").Value = -7.262399999999993
$ws.Range("D92").Value = -6.024
$ws.Range("D96").Value = -7.501800000000004
$ws.Range("D98").Value = -8.514800000000003
$ws.Range("A101").Value = -20.96619999999998
$ws.Range("D101").Value = -7.653900000000005
$ws.Range("A102").Value = -19.2697
$ws.Range("D102").Value = -8.223499999999996
